# 产品销售统计表 - 单肩包
# Add "最大销售利润" (max sales profit) and "最小销售利润" (min sales profit)
# summary cells next to the existing data table (columns I/J, rows 1-2),
# matching the currency formatting already used for column H (销售利润).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFormat = '"¥"#,##0.00;"¥"\-#,##0.00'

# Row 1: label + max value of the 销售利润 column (H2:H10)
$ws.Range("I1").Value = "最大销售利润"
$ws.Range("J1").Value = 51480
$ws.Range("J1").NumberFormat = $currencyFormat

# Row 2: label + min value of the 销售利润 column (H2:H10)
$ws.Range("I2").Value = "最小销售利润"
$ws.Range("J2").Value = 13200
$ws.Range("J2").NumberFormat = $currencyFormat

# Resize columns to fit the new content (as Excel does automatically).
$ws.Columns("A:J").AutoFit() | Out-Null
